$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Special characters used by these new "mex" coded-segment rows:
#  - Column A reuses the same bullet glyph already used throughout the
#    sheet (built from the real unicode char so the writer reuses the
#    existing shared string instead of minting a near-duplicate one).
#  - Some "Segment" (I) values embed literal control characters that
#    Excel's own _xHHHH_ escaping represents as e.g. "_x0003_256" in the
#    raw XML; typing that literal text would get double-escaped on save,
#    so the real control characters are built here instead.
$bullet = [string][char]0x25CF
$ctrl3  = [string][char]0x0003
$ctrl4  = [string][char]0x0004
$dagger = [string][char]0x2021

# New row data (columns A-M) for rows 208-215, mirroring the additional
# "mex" coded segments added in this run. B/C (Comment / Document group)
# are left blank like the neighboring rows, H (Weight score) is always 0,
# J/K hold the coverage count / coverage %, L is the coder name and M the
# creation timestamp. D is forced to text since these sheets always store
# the numeric-looking Document name as text.
$rows = @(
    @{ Row=208; D="2302"; F="3: 1686"; G="3: 1689"; I=($ctrl3 + "256"); J=4; K=0.01244516349833545837; M="1/31/19 13:54:31" },
    @{ Row=209; D="2628"; F="2: 6182"; G="2: 6185"; I=($ctrl4 + "256"); J=4; K=0.02462902530632350012; M="1/31/19 13:55:01" },
    @{ Row=210; D="2628"; F="2: 6202"; G="2: 6205"; I=($ctrl4 + "256"); J=4; K=0.02462902530632350012; M="1/31/19 13:55:07" },
    @{ Row=211; D="2628"; F="2: 6217"; G="2: 6220"; I=($ctrl4 + "256"); J=4; K=0.02462902530632350012; M="1/31/19 13:55:11" },
    @{ Row=212; D="2628"; F="2: 6254"; G="2: 6257"; I=($ctrl4 + "256"); J=4; K=0.02462902530632350012; M="1/31/19 13:55:15" },
    @{ Row=213; D="2628"; F="2: 6120"; G="2: 6122"; I=($ctrl4 + "32");  J=3; K=0.01847176897974262683; M="1/31/19 13:55:31" },
    @{ Row=214; D="2628"; F="2: 6135"; G="2: 6137"; I=($ctrl4 + "32");  J=3; K=0.01847176897974262683; M="1/31/19 13:55:35" },
    @{ Row=215; D="3872"; F="2: 3375"; G="2: 3377"; I=($dagger + "32"); J=3; K=0.0133838947133615889;  M="1/31/19 13:56:14" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Column D (Document name) looks numeric ("2302", "2628", ...) but must
    # be stored as text, like every other row in this sheet. Borrow a
    # text-formatted cell's number format (row 206, column C already uses
    # numFmtId 49 "@") before assigning the value so the digits aren't
    # auto-converted to a number - this avoids Excel minting a brand new
    # "quote prefixed" style the way a literal leading apostrophe would.
    $ws.Range("C206").Copy()
    $ws.Range("D" + $rowNum).PasteSpecial(-4122)
    $ws.Range("D" + $rowNum).Value = $r.D

    # Set the remaining cell contents (plain Value assignment), so each
    # cell lands with the right stored type (shared string vs. number).
    $ws.Range("A" + $rowNum).Value = $bullet
    $ws.Range("E" + $rowNum).Value = "MIC"
    $ws.Range("F" + $rowNum).Value = $r.F
    $ws.Range("G" + $rowNum).Value = $r.G
    $ws.Range("H" + $rowNum).Value = 0
    $ws.Range("I" + $rowNum).Value = $r.I
    $ws.Range("J" + $rowNum).Value = $r.J
    $ws.Range("K" + $rowNum).Value = $r.K
    $ws.Range("L" + $rowNum).Value = "chen"
    $ws.Range("M" + $rowNum).Value = $r.M

    # Now copy formatting (styles, fills, borders) from row 206, which
    # already carries the same per-column style pattern needed here. A
    # formats-only paste keeps the values/types set above intact.
    $ws.Range("A206:M206").Copy()
    $ws.Range("A" + $rowNum + ":M" + $rowNum).PasteSpecial(-4122)

    $ws.Range("A" + $rowNum).EntireRow.RowHeight = 16
}
